# Switch to use new prep = TRUE option, add EF table
#
# The "Slurry" sheet used to carry two identifying columns
# (man.name, man.source) ahead of the acid/man.dm/man.ph data. The new
# prep step drops the man.name column entirely: man.source becomes the
# first column, and every row now records the manure source as
# "Afgasset biomasse" (previously "Digestate").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# Remove the old "man.name" column (A). This shifts man.source, acid,
# man.dm and man.ph one column to the left, carrying their styles with
# them, and updates the sheet dimension from A1:E5 to A1:D5.
$ws.Columns("A").Delete()

# The values that used to live in man.source (now column A) said
# "Digestate" - update them to reflect the new source name.
$ws.Range("A2:A5").Value = "Afgasset biomasse"

# Match the author's last-saved selection on this sheet.
$ws.Range("E12").Select() | Out-Null
